$d = $word.ActiveDocument

$replacements = @(
    @("595×9=5355", "787×7=5509"),
    @("152×6=912", "346×2=692"),
    @("806×6=4836", "467×4=1868"),
    @("542×8=4336", "780×7=5460"),
    @("132×5=660", "208×6=1248"),
    @("294×4=1176", "910×6=5460"),
    @("857×4=3428", "317×9=2853"),
    @("161×5=805", "547×7=3829"),
    @("920×8=7360", "652×2=1304"),
    @("536×5=2680", "406×5=2030"),
    @("422×8=3376", "491×8=3928"),
    @("548×6=3288", "985×4=3940"),
    @("525×5=2625", "356×8=2848"),
    @("635×5=3175", "435×4=1740"),
    @("493×5=2465", "447×9=4023"),
    @("469×2=938", "506×8=4048"),
    @("185×3=555", "354×6=2124"),
    @("187×8=1496", "977×8=7816"),
    @("432×9=3888", "411×6=2466"),
    @("400×7=2800", "658×8=5264"),
    @("762×9=6858", "869×6=5214"),
    @("564×8=4512", "730×7=5110"),
    @("591×3=1773", "506×7=3542"),
    @("459×3=1377", "497×8=3976"),
    @("779×3=2337", "841×3=2523")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
